$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Seed the new "user:prod_*" attribute-name strings in the desired order before
# rewriting the column B cells below (keeps the shared-string table ordering stable)
$seedCell = $ws.Range("Z1")
$seedCell.Value = "user:prod_type"
$seedCell.Value = "user:prod_material"
$seedCell.Value = "user:prod_size"
$seedCell.Value = "user:prod_mat_attrib_1"
$seedCell.Value = "user:prod_model_type"
$seedCell.Value = "user:prod_material_attrib_2"
$seedCell.Value = "user:prod_manufacturer"
$seedCell.Value = "user:prod_design_type"
$seedCell.Value = "user:prod_subject"
$seedCell.Value = "user:prod_strategy"
$seedCell.Value = "user:prod_publisher"
$seedCell.ClearContents()

# Rename legacy "prod_*" filter attribute names to the new "user:prod_*" namespace
$ws.Cells.Item(2, 2).Value = "user:prod_type"
$ws.Cells.Item(3, 2).Value = "user:prod_material"
$ws.Cells.Item(4, 2).Value = "user:prod_size"
$ws.Cells.Item(5, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(6, 2).Value = "user:prod_material_attrib_2"
$ws.Cells.Item(7, 2).Value = "user:prod_model_type"
$ws.Cells.Item(9, 2).Value = "user:prod_type"
$ws.Cells.Item(10, 2).Value = "user:prod_material"
$ws.Cells.Item(11, 2).Value = "user:prod_size"
$ws.Cells.Item(12, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(13, 2).Value = "user:prod_design_type"
$ws.Cells.Item(14, 2).Value = "user:prod_model_type"
$ws.Cells.Item(16, 2).Value = "user:prod_type"
$ws.Cells.Item(17, 2).Value = "user:prod_material"
$ws.Cells.Item(18, 2).Value = "user:prod_model_type"
$ws.Cells.Item(19, 2).Value = "user:prod_size"
$ws.Cells.Item(20, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(22, 2).Value = "user:prod_type"
$ws.Cells.Item(23, 2).Value = "user:prod_material"
$ws.Cells.Item(24, 2).Value = "user:prod_size"
$ws.Cells.Item(25, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(26, 2).Value = "user:prod_material_attrib_2"
$ws.Cells.Item(28, 2).Value = "user:prod_type"
$ws.Cells.Item(29, 2).Value = "user:prod_material"
$ws.Cells.Item(30, 2).Value = "user:prod_model_type"
$ws.Cells.Item(31, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(32, 2).Value = "user:prod_manufacturer"
$ws.Cells.Item(34, 2).Value = "user:prod_type"
$ws.Cells.Item(35, 2).Value = "user:prod_material"
$ws.Cells.Item(36, 2).Value = "user:prod_design_type"
$ws.Cells.Item(37, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(38, 2).Value = "user:prod_size"
$ws.Cells.Item(40, 2).Value = "user:prod_type"
$ws.Cells.Item(41, 2).Value = "user:prod_material"
$ws.Cells.Item(42, 2).Value = "user:prod_size"
$ws.Cells.Item(43, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(45, 2).Value = "user:prod_type"
$ws.Cells.Item(46, 2).Value = "user:prod_material"
$ws.Cells.Item(47, 2).Value = "user:prod_design_type"
$ws.Cells.Item(48, 2).Value = "user:prod_size"
$ws.Cells.Item(49, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(51, 2).Value = "user:prod_type"
$ws.Cells.Item(52, 2).Value = "user:prod_material"
$ws.Cells.Item(53, 2).Value = "user:prod_size"
$ws.Cells.Item(54, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(56, 2).Value = "user:prod_type"
$ws.Cells.Item(57, 2).Value = "user:prod_material"
$ws.Cells.Item(58, 2).Value = "user:prod_model_type"
$ws.Cells.Item(59, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(60, 2).Value = "user:prod_manufacturer"
$ws.Cells.Item(62, 2).Value = "user:prod_type"
$ws.Cells.Item(63, 2).Value = "user:prod_material"
$ws.Cells.Item(64, 2).Value = "user:prod_size"
$ws.Cells.Item(65, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(66, 2).Value = "user:prod_manufacturer"
$ws.Cells.Item(68, 2).Value = "user:prod_type"
$ws.Cells.Item(69, 2).Value = "user:prod_material"
$ws.Cells.Item(70, 2).Value = "user:prod_size"
$ws.Cells.Item(71, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(72, 2).Value = "user:prod_manufacturer"
$ws.Cells.Item(74, 2).Value = "user:prod_type"
$ws.Cells.Item(75, 2).Value = "user:prod_material"
$ws.Cells.Item(76, 2).Value = "user:prod_design_type"
$ws.Cells.Item(77, 2).Value = "user:prod_size"
$ws.Cells.Item(78, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(80, 2).Value = "user:prod_type"
$ws.Cells.Item(81, 2).Value = "user:prod_material"
$ws.Cells.Item(82, 2).Value = "user:prod_design_type"
$ws.Cells.Item(83, 2).Value = "user:prod_material_attrib_2"
$ws.Cells.Item(84, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(86, 2).Value = "user:prod_type"
$ws.Cells.Item(87, 2).Value = "user:prod_material"
$ws.Cells.Item(88, 2).Value = "user:prod_size"
$ws.Cells.Item(89, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(91, 2).Value = "user:prod_type"
$ws.Cells.Item(92, 2).Value = "user:prod_material"
$ws.Cells.Item(93, 2).Value = "user:prod_design_type"
$ws.Cells.Item(94, 2).Value = "user:prod_size"
$ws.Cells.Item(95, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(97, 2).Value = "user:prod_type"
$ws.Cells.Item(98, 2).Value = "user:prod_material"
$ws.Cells.Item(99, 2).Value = "user:prod_design_type"
$ws.Cells.Item(100, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(102, 2).Value = "user:prod_type"
$ws.Cells.Item(103, 2).Value = "user:prod_material"
$ws.Cells.Item(104, 2).Value = "user:prod_size"
$ws.Cells.Item(105, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(106, 2).Value = "user:prod_material_attrib_2"
$ws.Cells.Item(108, 2).Value = "user:prod_type"
$ws.Cells.Item(109, 2).Value = "user:prod_material"
$ws.Cells.Item(110, 2).Value = "user:prod_size"
$ws.Cells.Item(111, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(112, 2).Value = "user:prod_material_attrib_2"
$ws.Cells.Item(114, 2).Value = "user:prod_type"
$ws.Cells.Item(115, 2).Value = "user:prod_material"
$ws.Cells.Item(116, 2).Value = "user:prod_design_type"
$ws.Cells.Item(117, 2).Value = "user:prod_material_attrib_2"
$ws.Cells.Item(118, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(120, 2).Value = "user:prod_type"
$ws.Cells.Item(121, 2).Value = "user:prod_material"
$ws.Cells.Item(122, 2).Value = "user:prod_design_type"
$ws.Cells.Item(123, 2).Value = "user:prod_mat_attrib_1"
$ws.Cells.Item(125, 2).Value = "user:prod_type"
$ws.Cells.Item(126, 2).Value = "user:prod_material"
$ws.Cells.Item(127, 2).Value = "user:prod_subject"
$ws.Cells.Item(128, 2).Value = "user:prod_strategy"
$ws.Cells.Item(129, 2).Value = "user:prod_publisher"
$ws.Cells.Item(131, 2).Value = "user:prod_type"
$ws.Cells.Item(132, 2).Value = "user:prod_material"
$ws.Cells.Item(133, 2).Value = "user:prod_subject"
$ws.Cells.Item(134, 2).Value = "user:prod_strategy"
$ws.Cells.Item(135, 2).Value = "user:prod_publisher"

# Update the saved selection/scroll position
$ws.Range("B129").Select()
